# Debit Note template: add a "Tipe" lookup sheet and a "Tipe" column (W)
# on Sheet1 with a data-validation dropdown sourced from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Add the "Tipe" worksheet after Sheet1 -----------------------------
$tipe = $wb.Worksheets.Add($null, $ws)
$tipe.Name = "Tipe"

$tipe.Range("B2").Value = "Listrik"
$tipe.Range("B3").Value = "Rent"
$tipe.Range("B4").Value = "Service"

# --- 2. Add the "Tipe" header + values on Sheet1, column W ----------------
$ws.Range("W4").Value = "Tipe"
$ws.Range("W4").Style = "Normal"
$ws.Range("W4").Interior.Color = 6299648
$ws.Range("W4").Font.Bold = $true
$ws.Range("W4").Font.Color = 16777215
$ws.Range("W4").HorizontalAlignment = -4108
$ws.Range("W4").VerticalAlignment = -4108
$ws.Range("W4").Borders.LineStyle = 1

$ws.Range("W5").Value = "Rent"
$ws.Range("W6").Value = "Listrik"
$ws.Range("W5:W6").Borders.LineStyle = 1

# --- 3. Data validation dropdown on W5:W6 ----------------------------------
$ws.Range("W5:W6").Validation.Add(3, 1, 1, "=Tipe!`$B`$2:`$B`$4")

# --- 4. Restore selection / active sheet -----------------------------------
$ws.Range("W5").Select()
$ws.Activate()
